$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: S.N. value changed from 9 to 12
$ws.Range("A13").Value = 12

# New row 14 appended to the change log
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 43664
$ws.Range("B14").NumberFormat = $ws.Range("B13").NumberFormat
$ws.Range("C14").Value = $ws.Range("C13").Text
$ws.Range("C14").NumberFormat = $ws.Range("C13").NumberFormat
$ws.Range("D14").Value = "Removed constraints for yield and edited version number to 30022019"
$ws.Range("E14").Value = $ws.Range("E13").Text
$ws.Range("F14").Value = "Kathmandu, Nepal"

# Update the selected cell to reflect the new last row
[void]$ws.Range("A15").Select()
